# "Generate Report for Archive"
# The localization status report is regenerated: every row whose status
# was "Ready for handoff" has now moved on to "In Translation". The
# "Ready for handoff" string stops being used anywhere in the workbook.
#
# Affected cells (discovered from the data, not hard-coded blindly):
#   Overview sheet : E5:F5, E6:F6, E7:F7   (zh-cn / de-de status columns)
#   zh-cn sheet    : C5, C6, C7            (Status column)
#   de-de sheet    : C5, C6, C7            (Status column)

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$oldStatus = "Ready for handoff"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ($cell.Value2 -eq $oldStatus) {
                $cell.Value2 = $newStatus
            }
        }
    }
    # Re-fit the status column(s) now that the text is shorter.
    $ws.Columns.AutoFit()
}
